# Auto-generated edit script applying the scheduled-runner value updates
# to the Chocobo_Profits crafting-profession sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Each sheet has columns H-N holding price/profit
# data that was refreshed by the runner; some rows gain or lose a cell
# (e.g. a previously-absent N value gets added, or a stale N value is
# cleared once M takes over) exactly as recorded upstream.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 919546
$ws.Range("I116").Value = 10000000
$ws.Range("J116").Value = 11500.6
$ws.Range("K116").Value = 10000000
$ws.Range("L116").Value = 11500.6
$ws.Range("M116").Value = -9996558
$ws.Range("N116").Value = -18384.6
$ws.Range("H138").Value = 1970.1177
$ws.Range("I138").Value = 1541
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 4623
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 517
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 29000
$ws.Range("J64").Value = 29000
$ws.Range("L64").Value = 29000
$ws.Range("N64").Value = -29496
$ws.Range("H67").Value = 29000
$ws.Range("J67").Value = 29000
$ws.Range("L67").Value = 29000
$ws.Range("N67").Value = -30716
$ws.Range("H123").Value = 48443
$ws.Range("J123").Value = 48443
$ws.Range("L123").Value = 48443
$ws.Range("N123").Value = -58243
$ws.Range("H131").Value = 42261.11
$ws.Range("J131").Value = 42261.11
$ws.Range("L131").Value = 42261.11
$ws.Range("N131").Value = -52341.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864
$ws.Range("H109").Value = 30000
$ws.Range("I109").Value = 30000
$ws.Range("K109").Value = 30000
$ws.Range("M109").Value = -28613
$ws.Range("H134").Value = 4190.8
$ws.Range("I134").Value = 1646
$ws.Range("J134").Value = 5039.067
$ws.Range("K134").Value = 4938
$ws.Range("L134").Value = 15117.201
$ws.Range("M134").Value = -2403
$ws.Range("N134").Value = -20187.201

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 184654.34
$ws.Range("I31").Value = 347384
$ws.Range("J31").Value = 3327.0286
$ws.Range("K31").Value = 347384
$ws.Range("L31").Value = 3327.0286
$ws.Range("M31").Value = -347089
$ws.Range("N31").Value = -3917.0286
$ws.Range("H34").Value = 184654.34
$ws.Range("I34").Value = 347384
$ws.Range("J34").Value = 3327.0286
$ws.Range("K34").Value = 347384
$ws.Range("L34").Value = 3327.0286
$ws.Range("M34").Value = -347182
$ws.Range("N34").Value = -3731.0286
$ws.Range("H86").Value = 2450.875
$ws.Range("I86").Value = 2038.1818
$ws.Range("J86").Value = 3358.8
$ws.Range("K86").Value = 2038.1818
$ws.Range("L86").Value = 3358.8
$ws.Range("M86").Value = -915.1818000000001
$ws.Range("N86").Value = -5604.8
$ws.Range("H89").Value = 2450.875
$ws.Range("I89").Value = 2038.1818
$ws.Range("J89").Value = 3358.8
$ws.Range("K89").Value = 10190.909
$ws.Range("L89").Value = 16794
$ws.Range("M89").Value = -4574.909
$ws.Range("N89").Value = -28026
$ws.Range("H132").Value = 3905.6
$ws.Range("I132").Value = 2885.6843
$ws.Range("K132").Value = 8657.052899999999
$ws.Range("M132").Value = -6127.052899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 464.59573
$ws.Range("I113").Value = 465.7931
$ws.Range("J113").Value = 462.66666
$ws.Range("K113").Value = 1397.3793
$ws.Range("L113").Value = 1387.99998
$ws.Range("M113").Value = 772.6206999999999
$ws.Range("N113").Value = -5727.999980000001
$ws.Range("H131").Value = 773.8817
$ws.Range("I131").Value = 407.55554
$ws.Range("J131").Value = 813.131
$ws.Range("K131").Value = 1222.66662
$ws.Range("L131").Value = 2439.393
$ws.Range("M131").Value = 3817.33338
$ws.Range("N131").Value = -12519.393

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4409.7144
$ws.Range("I40").Value = 3662.375
$ws.Range("J40").Value = 6040.273
$ws.Range("K40").Value = 3662.375
$ws.Range("L40").Value = 6040.273
$ws.Range("M40").Value = -3526.375
$ws.Range("N40").Value = -6312.273
$ws.Range("H109").Value = 29400
$ws.Range("J109").Value = 29400
$ws.Range("L109").Value = 29400
$ws.Range("N109").Value = -32174
$ws.Range("H122").Value = 3775.348
$ws.Range("I122").Value = 2768.2222
$ws.Range("J122").Value = 7401
$ws.Range("K122").Value = 8304.6666
$ws.Range("L122").Value = 22203
$ws.Range("M122").Value = -5854.6666
$ws.Range("N122").Value = -27103

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 18199.916
$ws.Range("J22").Value = 18199.916
$ws.Range("L22").Value = 18199.916
$ws.Range("N22").Value = -18785.916
$ws.Range("H69").Value = 2623
$ws.Range("I69").Value = 2623
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 2623
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -1874
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 2623
$ws.Range("I72").Value = 2623
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 7869
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -4125
$ws.Range("N72").ClearContents()
$ws.Range("H81").Value = 1565.5385
$ws.Range("I81").Value = 1137.5
$ws.Range("K81").Value = 2275
$ws.Range("M81").Value = -1214
$ws.Range("H84").Value = 1565.5385
$ws.Range("I84").Value = 1137.5
$ws.Range("K84").Value = 11375
$ws.Range("M84").Value = -6071
$ws.Range("H96").Value = 93074504
$ws.Range("I96").Value = 125000776
$ws.Range("J96").Value = 7937776
$ws.Range("K96").Value = 125000776
$ws.Range("L96").Value = 7937776
$ws.Range("M96").Value = -124999403
$ws.Range("N96").Value = -7940522
$ws.Range("H107").Value = 573.05554
$ws.Range("I107").Value = 398.95834
$ws.Range("J107").Value = 921.25
$ws.Range("K107").Value = 1196.87502
$ws.Range("L107").Value = 2763.75
$ws.Range("M107").Value = 723.1249800000001
$ws.Range("N107").Value = -6603.75
$ws.Range("H123").Value = 30631.062
$ws.Range("J123").Value = 30631.062
$ws.Range("L123").Value = 30631.062
$ws.Range("N123").Value = -40431.06200000001
$ws.Range("H132").Value = 8132501.5
$ws.Range("I132").Value = 990.7646999999999
$ws.Range("J132").Value = 13892322
$ws.Range("K132").Value = 2972.2941
$ws.Range("L132").Value = 41676966
$ws.Range("M132").Value = -442.2941000000001
$ws.Range("N132").Value = -41682026
$ws.Range("H136").Value = 3201.4814
$ws.Range("I136").Value = 682
$ws.Range("J136").Value = 6350.8335
$ws.Range("K136").Value = 2046
$ws.Range("L136").Value = 19052.5005
$ws.Range("M136").Value = 504
$ws.Range("N136").Value = -24152.5005
